$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7998.3335
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H132").Value = 2355399.2
$ws.Range("I132").Value = 2477204.5
$ws.Range("K132").Value = 7431613.5
$ws.Range("M132").Value = -7429083.5
$ws.Range("H138").Value = 1830.2449
$ws.Range("I138").Value = 1339.6842
$ws.Range("J138").Value = 3524.9092
$ws.Range("K138").Value = 4019.0526
$ws.Range("L138").Value = 10574.7276
$ws.Range("M138").Value = 1120.9474
$ws.Range("N138").Value = -20854.7276
$ws.Range("H141").Value = 2215.162
$ws.Range("I141").Value = 1815.742
$ws.Range("K141").Value = 5447.226
$ws.Range("M141").Value = -267.2259999999997
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2477.4211
$ws.Range("I2").Value = 2042.2
$ws.Range("K2").Value = 2042.2
$ws.Range("M2").Value = -1929.2
$ws.Range("H45").Value = 4037.6924
$ws.Range("I45").Value = 2910.8572
$ws.Range("K45").Value = 2910.8572
$ws.Range("M45").Value = -2533.8572
$ws.Range("H52").Value = 74999.5
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 74999.5
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 74999.5
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = -75635.5
$ws.Range("H61").Value = 5416.1787
$ws.Range("I61").Value = 1140.2
$ws.Range("J61").Value = 16106.125
$ws.Range("K61").Value = 1140.2
$ws.Range("L61").Value = 16106.125
$ws.Range("M61").Value = -928.2
$ws.Range("N61").Value = -16530.125
$ws.Range("H74").Value = 146302.27
$ws.Range("I74").Value = 167686.58
$ws.Range("K74").Value = 167686.58
$ws.Range("M74").Value = -166812.58
$ws.Range("H77").Value = 146302.27
$ws.Range("I77").Value = 167686.58
$ws.Range("K77").Value = 838432.8999999999
$ws.Range("M77").Value = -834064.8999999999
$ws.Range("H116").Value = 2477.4211
$ws.Range("I116").Value = 2042.2
$ws.Range("K116").Value = 2042.2
$ws.Range("M116").Value = 251.8
$ws.Range("H132").Value = 1307.8588
$ws.Range("I132").Value = 996.3099
$ws.Range("K132").Value = 2988.9297
$ws.Range("M132").Value = -458.9296999999997
$ws.Range("H136").Value = 5416.1787
$ws.Range("I136").Value = 1140.2
$ws.Range("J136").Value = 16106.125
$ws.Range("K136").Value = 3420.6
$ws.Range("L136").Value = 48318.375
$ws.Range("M136").Value = -870.6000000000004
$ws.Range("N136").Value = -53418.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2477.4211
$ws.Range("I3").Value = 2042.2
$ws.Range("K3").Value = 2042.2
$ws.Range("M3").Value = -1928.2
$ws.Range("H94").Value = 2542.7144
$ws.Range("I94").Value = 1183.1666
$ws.Range("K94").Value = 1183.1666
$ws.Range("M94").Value = -732.1666
$ws.Range("H105").Value = 3710.9722
$ws.Range("I105").Value = 3890.2
$ws.Range("K105").Value = 3890.2
$ws.Range("M105").Value = -2143.2
$ws.Range("H107").Value = 18277.516
$ws.Range("I107").Value = 26427.6
$ws.Range("K107").Value = 26427.6
$ws.Range("M107").Value = -24507.6
$ws.Range("H134").Value = 1561.375
$ws.Range("I134").Value = 1393.6818
$ws.Range("K134").Value = 4181.0454
$ws.Range("M134").Value = -1646.0454
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 850.15
$ws.Range("J58").Value = 722
$ws.Range("L58").Value = 722
$ws.Range("N58").Value = -1128
$ws.Range("H99").Value = 4845.4
$ws.Range("I99").Value = 4047.75
$ws.Range("J99").Value = 8036
$ws.Range("K99").Value = 4047.75
$ws.Range("L99").Value = 8036
$ws.Range("M99").Value = -2549.75
$ws.Range("N99").Value = -11032
$ws.Range("H107").Value = 2661.0908
$ws.Range("I107").Value = 275.6
$ws.Range("J107").Value = 4649
$ws.Range("K107").Value = 275.6
$ws.Range("L107").Value = 4649
$ws.Range("M107").Value = 1644.4
$ws.Range("N107").Value = -8489
$ws.Range("H126").Value = 4845.4
$ws.Range("I126").Value = 4047.75
$ws.Range("J126").Value = 8036
$ws.Range("K126").Value = 12143.25
$ws.Range("L126").Value = 24108
$ws.Range("M126").Value = -9673.25
$ws.Range("N126").Value = -29048
$ws.Range("H132").Value = 24821.77
$ws.Range("I132").Value = 29843.143
$ws.Range("K132").Value = 89529.429
$ws.Range("M132").Value = -86999.429
$ws.Range("H134").Value = 1211.9814
$ws.Range("I134").Value = 932.3958
$ws.Range("K134").Value = 2797.1874
$ws.Range("M134").Value = -262.1873999999998
$ws.Range("H136").Value = 850.15
$ws.Range("J136").Value = 722
$ws.Range("L136").Value = 2166
$ws.Range("N136").Value = -7266
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 3799.3333
$ws.Range("I6").Value = 4499.2
$ws.Range("K6").Value = 13497.6
$ws.Range("M6").Value = -13384.6
$ws.Range("H117").Value = 535.8333
$ws.Range("I117").Value = 535.8333
$ws.Range("K117").Value = 1607.4999
$ws.Range("M117").Value = 1834.5001
$ws.Range("H121").Value = 53739.81
$ws.Range("I121").Value = 101256.18
$ws.Range("K121").Value = 303768.54
$ws.Range("M121").Value = -302458.54
$ws.Range("H131").Value = 123323.83
$ws.Range("J131").Value = 1976.4482
$ws.Range("L131").Value = 5929.3446
$ws.Range("N131").Value = -16009.3446
$ws.Range("H132").Value = 2079.2
$ws.Range("I132").Value = 3249.5
$ws.Range("J132").Value = 1653.6364
$ws.Range("K132").Value = 29245.5
$ws.Range("L132").Value = 14882.7276
$ws.Range("M132").Value = -26715.5
$ws.Range("N132").Value = -19942.7276
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23569.857
$ws.Range("J15").Value = 23569.857
$ws.Range("L15").Value = 23569.857
$ws.Range("N15").Value = -24145.857
$ws.Range("H80").Value = 4114.148
$ws.Range("I80").Value = 2272.3157
$ws.Range("K80").Value = 2272.3157
$ws.Range("M80").Value = -1274.3157
$ws.Range("H81").Value = 23569.857
$ws.Range("J81").Value = 23569.857
$ws.Range("L81").Value = 23569.857
$ws.Range("N81").Value = -25565.857
$ws.Range("H83").Value = 4114.148
$ws.Range("I83").Value = 2272.3157
$ws.Range("K83").Value = 11361.5785
$ws.Range("M83").Value = -6369.5785
$ws.Range("H84").Value = 23569.857
$ws.Range("J84").Value = 23569.857
$ws.Range("L84").Value = 70709.571
$ws.Range("N84").Value = -80693.571
$ws.Range("H132").Value = 2519.5625
$ws.Range("I132").Value = 2096.2593
$ws.Range("J132").Value = 4805.4
$ws.Range("K132").Value = 6288.777900000001
$ws.Range("L132").Value = 14416.2
$ws.Range("M132").Value = -3758.777900000001
$ws.Range("N132").Value = -19476.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1128.5
$ws.Range("J22").Value = 1455
$ws.Range("L22").Value = 1455
$ws.Range("N22").Value = -2045
$ws.Range("H27").Value = 1128.5
$ws.Range("J27").Value = 1455
$ws.Range("L27").Value = 1455
$ws.Range("N27").Value = -1669
$ws.Range("H46").Value = 5446.4736
$ws.Range("I46").Value = 2740.4
$ws.Range("J46").Value = 6412.9287
$ws.Range("K46").Value = 2740.4
$ws.Range("L46").Value = 6412.9287
$ws.Range("M46").Value = -2552.4
$ws.Range("N46").Value = -6788.9287
$ws.Range("H132").Value = 3255.818
$ws.Range("J132").Value = 4998.5
$ws.Range("L132").Value = 14995.5
$ws.Range("N132").Value = -20055.5
$ws.Range("H136").Value = 3024.724
$ws.Range("I136").Value = 2652.157
$ws.Range("K136").Value = 7956.471
$ws.Range("M136").Value = -5406.471
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 28999
$ws.Range("I40").Value = 28999
$ws.Range("K40").Value = 28999
$ws.Range("M40").Value = -28850
$ws.Range("H62").Value = 4666.3335
$ws.Range("I62").Value = 3999.5
$ws.Range("K62").Value = 3999.5
$ws.Range("M62").Value = -3375.5
$ws.Range("H65").Value = 4666.3335
$ws.Range("I65").Value = 3999.5
$ws.Range("K65").Value = 19997.5
$ws.Range("M65").Value = -16877.5
$ws.Range("H107").Value = 1433.6
$ws.Range("I107").Value = 1167
$ws.Range("K107").Value = 3501
$ws.Range("M107").Value = -1581
$ws.Range("H122").Value = 15192504
$ws.Range("I122").Value = 15666644
$ws.Range("K122").Value = 46999932
$ws.Range("M122").Value = -46997482
$ws.Range("H132").Value = 5119711.5
$ws.Range("I132").Value = 8360565.5
$ws.Range("K132").Value = 25081696.5
$ws.Range("M132").Value = -25079166.5
$ws.Range("H136").Value = 9564.245999999999
$ws.Range("I136").Value = 10348.5
$ws.Range("K136").Value = 31045.5
$ws.Range("M136").Value = -28495.5
